# Apply the updated imputed values for column A (and some column E) as per
# the "Update Name of Algo" commit (KNN re-run produced slightly different
# imputed values for a handful of cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4"   = -20.229
    "A7"   = -20.053
    "E7"   = 16.327
    "E15"  = 16.113
    "A16"  = -22.033
    "E21"  = 16.414
    "E22"  = 16.602
    "E23"  = 16.434
    "A28"  = -22.087
    "A29"  = -21.344
    "A32"  = -21.816
    "E34"  = 16.805
    "A40"  = -19.965
    "E43"  = 17.083
    "E45"  = 16.774
    "E50"  = 16.403
    "E51"  = 16.828
    "A52"  = -21.963
    "A57"  = -22.275
    "A66"  = -21.546
    "E66"  = 17.394
    "E67"  = 17.325
    "E79"  = 17.07
    "E84"  = 16.608
    "E92"  = 17.885
    "E97"  = 16.769
    "A100" = -22.121
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$wb.Save()
